# Apply the "First selectable test ready" edit to UserData.xlsx
#
# Summary of the change:
#  - Rename the Droppable-test header columns D1:F1 from ItemCat1/2/3 to
#    Item1/2/3 on the "DataSetInteractionPages" sheet.
#  - Give the whole data block (B2:F14) a Text ("@") number format, which
#    is what stamps style index 2 onto the existing cells and materialises
#    the new (blank) spacer rows 7, 9, 12 and 14 with that same style.
#  - Update the Selectable test row: rename the key from "...SelectThree..."
#    to "...SelectTwo...", and replace the old "Item 5"/"Item 2"/"Item 1"
#    text values in D13:F13 with the numeric list positions 1 and 5 (F13
#    is cleared out).
#  - Move the active/selected worksheet from "DataSetInteractionPages" to
#    "DataSetRegistrationUser".

$wb = $excel.ActiveWorkbook

$wsData        = $wb.Worksheets.Item("DataSetInteractionPages")
$wsRegistration = $wb.Worksheets.Item("DataSetRegistrationUser")

# --- Selectable test row (13): new key name + numeric list values ---------
$wsData.Range("A13").Value = "SelectableItems_SelectTwo_SelectedElementsStatusChangedToSelected"

# --- Header row: ItemCat1/2/3 -> Item1/2/3 ---------------------------------
$wsData.Range("D1").Value = "Item1"
$wsData.Range("E1").Value = "Item2"
$wsData.Range("F1").Value = "Item3"

$wsData.Range("D13").Value = 1
$wsData.Range("E13").Value = 5
$wsData.Range("F13").ClearContents()

# --- Stamp the Text number format across the data block --------------------
# This both applies style 2 to the existing B:C value cells / D:F category
# cells and creates the new blank (but styled) spacer rows 7, 9, 12, 14.
$wsData.Range("B2:F14").NumberFormat = "@"

# --- Leave the cursor where the author left it on the data sheet, then -----
# --- switch the active sheet to DataSetRegistrationUser --------------------
[void]$wsData.Range("D19").Select()
$wsRegistration.Activate()
